# Power_Storage.xlsx
#
# 1) Recompute the "true" Excel column width for the handful of columns
#    whose stored <col width="..."> wasn't an exact multiple of 1/256 of a
#    character (i.e. wasn't a value Excel itself could ever have written):
#    col A (1), cols Y:Z (25:26) and cols AC:AD (29:30). Excel's COM model
#    exposes width as characters-of-the-Normal-style-font *after* removing
#    the fixed 5-pixel cell padding, so ColumnWidth = storedWidth - 5/6
#    (5/6 char ~= 5px at the default Calibri-11 digit width) round-trips
#    back to the corrected stored width on save.
# 2) Row 1 (the title row) now carries an explicit 24pt custom row height.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 5.5703125 - 5/6
$ws.Columns.Item(25).ColumnWidth = 20.140625 - 5/6
$ws.Columns.Item(26).ColumnWidth = 20.140625 - 5/6
$ws.Columns.Item(29).ColumnWidth = 24.5703125 - 5/6
$ws.Columns.Item(30).ColumnWidth = 24.5703125 - 5/6

$ws.Rows.Item(1).RowHeight = 24
